$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Remove the data rows on Sheet1, and clear the header row's values (keep formatting)
$ws1.Range("A2:F3").EntireRow.Delete()
$ws1.Range("A1:F1").ClearContents()

# Make Sheet1 the active sheet/tab with a specific selection
$ws1.Activate()
$ws1.Range("A1:G4").Select()
$ws1.Range("G4").Activate()
